$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so purely numeric-looking
# values (e.g. "218.43") are written back as text, matching the original
# inlineStr cell type instead of being coerced to a float by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.094.79"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.647.33"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "218.43"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "0.5199"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.2618"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "0.06297"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "20.29"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "0.07653"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "1.647.98"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "1.874.22"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "0.5568"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "0.0₅8111"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "65.07"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "26.053.36"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "4.585"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "194.08"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "5.916"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "145.11"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "0.05440"
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").Value = "1.268"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "3.433"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "3.319"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D36").Value = "2.780"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "0.9417"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "0.5583"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").Value = "0.01569"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "5.732"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("D42").Value = "1.027.45"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").Value = "0.8210"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "100.49"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "1.785.49"
$ws.Range("E46").Value = "  +9.20%  "
$ws.Range("D47").Value = "57.16"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "0.9987"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "0.4318"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "7.883"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "0.05099"
$ws.Range("E51").Value = "  -4.28%  "

# Restore the original (default) cell style now that the values are
# committed as text, so no visible/number-format change persists.
$ws.Range("D2:D51").Style = "Normal"
